# Auto-generated edit script: updates F/G numeric columns per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 339
$ws.Range("F3").Value = 226
$ws.Range("F4").Value = 553
$ws.Range("F5").Value = 1334
$ws.Range("F6").Value = 651
$ws.Range("F7").Value = 345
$ws.Range("F8").Value = 27
$ws.Range("F10").Value = 410
$ws.Range("F11").Value = 6174
$ws.Range("F12").Value = 113
$ws.Range("F14").Value = 1892
$ws.Range("F15").Value = 4628
$ws.Range("F19").Value = 5397
$ws.Range("F20").Value = 7043
$ws.Range("F22").Value = 1083
$ws.Range("F23").Value = 749
$ws.Range("F24").Value = 3971
$ws.Range("F25").Value = 546
$ws.Range("F26").Value = 73
$ws.Range("F27").Value = 227
$ws.Range("F29").Value = 1048
$ws.Range("F30").Value = 1486
$ws.Range("F31").Value = 547
$ws.Range("F32").Value = 677
$ws.Range("F33").Value = 1676
$ws.Range("F34").Value = 236
$ws.Range("F35").Value = 1867
$ws.Range("F37").Value = 1226
$ws.Range("F39").Value = 1335
$ws.Range("F40").Value = 677
$ws.Range("F41").Value = 317
$ws.Range("F42").Value = 724
$ws.Range("F43").Value = 3639
$ws.Range("F47").Value = 21
$ws.Range("F49").Value = 3946

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G9").Value = 228

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4354

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4354
$ws.Range("F3").Value = 339
$ws.Range("F7").Value = 226
$ws.Range("F8").Value = 553
$ws.Range("F10").Value = 1334
$ws.Range("G11").Value = 228
$ws.Range("F12").Value = 651
$ws.Range("F13").Value = 345
$ws.Range("F14").Value = 27
$ws.Range("F16").Value = 410
$ws.Range("F17").Value = 113
$ws.Range("F19").Value = 4628
$ws.Range("F20").Value = 5397
$ws.Range("F22").Value = 1083
$ws.Range("F23").Value = 749
$ws.Range("F24").Value = 3971
$ws.Range("F25").Value = 546
$ws.Range("F26").Value = 227
$ws.Range("F29").Value = 1048
$ws.Range("F30").Value = 1486
$ws.Range("F31").Value = 547
$ws.Range("F32").Value = 677
$ws.Range("F33").Value = 1676
$ws.Range("F34").Value = 1867
$ws.Range("F37").Value = 677
$ws.Range("F39").Value = 317
$ws.Range("F41").Value = 3639
$ws.Range("F46").Value = 21
$ws.Range("F49").Value = 3946
